$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1969.7941
$ws.Range("I17").Value = 725.6
$ws.Range("J17").Value = 2184.3103
$ws.Range("K17").Value = 2176.8
$ws.Range("L17").Value = 6552.9309
$ws.Range("M17").Value = -2008.8
$ws.Range("N17").Value = -6888.9309
$ws.Range("H111").Value = 54516.65
$ws.Range("I111").Value = 80727
$ws.Range("J111").Value = 5840.2856
$ws.Range("K111").Value = 242181
$ws.Range("L111").Value = 17520.8568
$ws.Range("M111").Value = -239114
$ws.Range("N111").Value = -23654.8568
$ws.Range("H116").Value = 9497.842000000001
$ws.Range("I116").Value = 4134
$ws.Range("J116").Value = 13398.818
$ws.Range("K116").Value = 4134
$ws.Range("L116").Value = 13398.818
$ws.Range("M116").Value = -692
$ws.Range("N116").Value = -20282.818
$ws.Range("H127").Value = 2414.5454
$ws.Range("I127").Value = 395.66666
$ws.Range("K127").Value = 1186.99998
$ws.Range("M127").Value = 3773.00002
$ws.Range("H137").Value = 3663.8333
$ws.Range("I137").Value = 3729.3
$ws.Range("J137").Value = 3582
$ws.Range("K137").Value = 11187.9
$ws.Range("L137").Value = 10746
$ws.Range("M137").Value = -8637.900000000001
$ws.Range("N137").Value = -15846
$ws.Range("H138").Value = 4409.117
$ws.Range("I138").Value = 1571.52
$ws.Range("K138").Value = 4714.559999999999
$ws.Range("M138").Value = 425.4400000000005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4908.4463
$ws.Range("I32").Value = 4908.4463
$ws.Range("K32").Value = 4908.4463
$ws.Range("M32").Value = -4621.4463
$ws.Range("H45").Value = 2493.3076
$ws.Range("I45").Value = 1823.8889
$ws.Range("K45").Value = 1823.8889
$ws.Range("M45").Value = -1446.8889
$ws.Range("H74").Value = 1181.2727
$ws.Range("I74").Value = 1169.45
$ws.Range("K74").Value = 1169.45
$ws.Range("M74").Value = -295.45
$ws.Range("H77").Value = 1181.2727
$ws.Range("I77").Value = 1169.45
$ws.Range("K77").Value = 5847.25
$ws.Range("M77").Value = -1479.25
$ws.Range("H122").Value = 5654.5625
$ws.Range("I122").Value = 4011.5
$ws.Range("J122").Value = 5889.2856
$ws.Range("K122").Value = 12034.5
$ws.Range("L122").Value = 17667.8568
$ws.Range("M122").Value = -9584.5
$ws.Range("N122").Value = -22567.8568
$ws.Range("H132").Value = 3369.85
$ws.Range("I132").Value = 3945.182
$ws.Range("K132").Value = 11835.546
$ws.Range("M132").Value = -9305.545999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 87214.25
$ws.Range("J59").Value = 87214.25
$ws.Range("L59").Value = 87214.25
$ws.Range("N59").Value = -88908.25
$ws.Range("H134").Value = 28892.025
$ws.Range("I134").Value = 1923.4412
$ws.Range("K134").Value = 5770.3236
$ws.Range("M134").Value = -3235.3236

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 27685
$ws.Range("J41").Value = 34240.11
$ws.Range("L41").Value = 34240.11
$ws.Range("N41").Value = -35096.11
$ws.Range("H99").Value = 5689.9443
$ws.Range("I99").Value = 3703.6667
$ws.Range("K99").Value = 3703.6667
$ws.Range("M99").Value = -2205.6667
$ws.Range("H126").Value = 5689.9443
$ws.Range("I126").Value = 3703.6667
$ws.Range("K126").Value = 11111.0001
$ws.Range("M126").Value = -8641.000100000001
$ws.Range("H132").Value = 2490.6316
$ws.Range("I132").Value = 1426.5714
$ws.Range("J132").Value = 5470
$ws.Range("K132").Value = 4279.7142
$ws.Range("L132").Value = 16410
$ws.Range("M132").Value = -1749.7142
$ws.Range("N132").Value = -21470
$ws.Range("H134").Value = 360559.72
$ws.Range("I134").Value = 3298.9167
$ws.Range("J134").Value = 2504124.5
$ws.Range("K134").Value = 9896.750100000001
$ws.Range("L134").Value = 7512373.5
$ws.Range("M134").Value = -7361.750100000001
$ws.Range("N134").Value = -7517443.5
$ws.Range("H139").Value = 99750
$ws.Range("J139").Value = 99750
$ws.Range("L139").Value = 99750
$ws.Range("N139").Value = -110030

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1873383
$ws.Range("I4").Value = 2385145.8
$ws.Range("J4").Value = 210154
$ws.Range("K4").Value = 7155437.399999999
$ws.Range("L4").Value = 630462
$ws.Range("M4").Value = -7155325.399999999
$ws.Range("N4").Value = -630686
$ws.Range("H131").Value = 3221.6924
$ws.Range("I131").Value = 1000
$ws.Range("J131").Value = 3265.255
$ws.Range("K131").Value = 3000
$ws.Range("L131").Value = 9795.764999999999
$ws.Range("M131").Value = 2040
$ws.Range("N131").Value = -19875.765

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1561.25
$ws.Range("I102").Value = 601.88
$ws.Range("K102").Value = 601.88
$ws.Range("M102").Value = 1020.12
$ws.Range("H122").Value = 5315.3
$ws.Range("I122").Value = 1997.5
$ws.Range("K122").Value = 5992.5
$ws.Range("M122").Value = -3542.5
$ws.Range("H126").Value = 4183
$ws.Range("I126").Value = 3325
$ws.Range("J126").Value = 4428.143
$ws.Range("K126").Value = 9975
$ws.Range("L126").Value = 13284.429
$ws.Range("M126").Value = -7505
$ws.Range("N126").Value = -18224.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 840701.25
$ws.Range("J7").Value = 1116601.8
$ws.Range("L7").Value = 1116601.8
$ws.Range("N7").Value = -1116825.8
$ws.Range("H40").Value = 5559477
$ws.Range("I40").Value = 6253786.5
$ws.Range("K40").Value = 6253786.5
$ws.Range("M40").Value = -6253650.5
$ws.Range("H48").Value = 22250
$ws.Range("I48").Value = 4500
$ws.Range("K48").Value = 4500
$ws.Range("M48").Value = -3839
$ws.Range("H54").Value = 39900
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 39900
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 39900
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -41188
$ws.Range("H82").Value = 633.3333
$ws.Range("I82").Value = 650
$ws.Range("J82").Value = 600
$ws.Range("K82").Value = 650
$ws.Range("L82").Value = 600
$ws.Range("M82").Value = -289
$ws.Range("N82").Value = -1322
$ws.Range("H85").Value = 633.3333
$ws.Range("I85").Value = 650
$ws.Range("J85").Value = 600
$ws.Range("K85").Value = 650
$ws.Range("L85").Value = 600
$ws.Range("M85").Value = 598
$ws.Range("N85").Value = -3096
$ws.Range("H93").Value = 2242.8333
$ws.Range("I93").Value = 2190.9092
$ws.Range("J93").Value = 2324.4285
$ws.Range("K93").Value = 2190.9092
$ws.Range("L93").Value = 2324.4285
$ws.Range("M93").Value = -942.9092000000001
$ws.Range("N93").Value = -4820.4285
$ws.Range("H126").Value = 840701.25
$ws.Range("J126").Value = 1116601.8
$ws.Range("L126").Value = 3349805.4
$ws.Range("N126").Value = -3354745.4
$ws.Range("H133").Value = 62842.715
$ws.Range("J133").Value = 64983.168
$ws.Range("L133").Value = 64983.168
$ws.Range("N133").Value = -70043.16800000001
$ws.Range("H136").Value = 2006778.8
$ws.Range("I136").Value = 2862054.2
$ws.Range("K136").Value = 8586162.600000001
$ws.Range("M136").Value = -8583612.600000001
$ws.Range("H140").Value = 75000
$ws.Range("J140").Value = 75000
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 6926.8335
$ws.Range("I74").Value = 4797.3335
$ws.Range("J74").Value = 7636.6665
$ws.Range("K74").Value = 4797.3335
$ws.Range("L74").Value = 7636.6665
$ws.Range("M74").Value = -3861.3335
$ws.Range("N74").Value = -9508.666499999999
$ws.Range("H77").Value = 6926.8335
$ws.Range("I77").Value = 4797.3335
$ws.Range("J77").Value = 7636.6665
$ws.Range("K77").Value = 14392.0005
$ws.Range("L77").Value = 22909.9995
$ws.Range("M77").Value = -9712.000499999998
$ws.Range("N77").Value = -32269.9995
$ws.Range("H132").Value = 17761.637
$ws.Range("I132").Value = 2536.3171
$ws.Range("J132").Value = 42731.16
$ws.Range("K132").Value = 7608.951300000001
$ws.Range("L132").Value = 128193.48
$ws.Range("M132").Value = -5078.951300000001
$ws.Range("N132").Value = -133253.48
$ws.Range("H138").Value = 83088.5
$ws.Range("J138").Value = 83088.5
$ws.Range("L138").Value = 83088.5
$ws.Range("N138").Value = -93368.5
